# ============================================================
# Add a new "2022-Q3" worksheet (with fund holdings detail)
# right after "总计", shifting all the quarter sheets down.
# Also insert the corresponding summary row in "总计".
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Locate source sheets by name --------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

# ---- 2. Duplicate the "2022-Q2" sheet (same column layout) right
#         after "总计" - this becomes our new "2022-Q3" sheet, already
#         carrying the correct header / styles.
$q2Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

$rows2022Q3 = @(
  @(0, "011855", "银华长荣混合", "10.55", "65.92", "6.59", "0.6952", 5),
  @(1, "013247", "交银瑞卓三年持有期混合", "19.44", "67.60", "2.16", "0.4199", 10),
  @(2, "470007", "汇添富上证综合指数", "7.41", "94.17", "2.96", "0.2193", 3),
  @(3, "510210", "富国上证综指ETF", "9.02", "99.24", "2.08", "0.1876", 5),
  @(4, "008978", "银华长丰混合", "2.24", "73.86", "5.93", "0.1328", 4),
  @(5, "008261", "招商研究优选股票A", "1.95", "87.39", "6.11", "0.1191", 5),
  @(6, "005706", "兴业龙腾双益平衡混合", "1.81", "32.02", "5.12", "0.0927", 3),
  @(7, "165310", "建信沪深300指数增强（LOF）A", "2.78", "93.07", "2.85", "0.0792", 4),
  @(8, "510760", "国泰上证综合ETF", "2.75", "94.96", "2.43", "0.0668", 2),
  @(9, "012877", "富荣福耀混合C", "1.77", "61.19", "3.49", "0.0618", 8),
  @(10, "012708", "东方红中证东方红红利低波动指数A", "3.27", "93.80", "1.62", "0.0530", 5),
  @(11, "013611", "工银民瑞一年持有混合A", "3.05", "21.97", "1.67", "0.0509", 1),
  @(12, "515300", "嘉实沪深300红利低波动ETF", "0.94", "99.19", "4.87", "0.0458", 1),
  @(13, "519677", "银河定投宝腾讯济安指数", "2.88", "92.40", "1.48", "0.0426", 6),
  @(14, "008262", "招商研究优选股票C", "0.58", "87.39", "6.11", "0.0354", 5),
  @(15, "011376", "华宝安享混合", "6.06", "20.12", "0.51", "0.0309", 7),
  @(16, "003154", "华宝新活力灵活配置混合", "4.40", "26.55", "0.60", "0.0264", 5),
  @(17, "005381", "泰康睿利量化多策略混合A", "0.48", "78.29", "4.92", "0.0236", 1),
  @(18, "005382", "泰康睿利量化多策略混合C", "0.48", "78.29", "4.92", "0.0236", 1),
  @(19, "660006", "农银大盘蓝筹混合", "1.27", "84.42", "1.60", "0.0203", 3),
  @(20, "003144", "华宝新机遇灵活配置混合（LOF）C", "4.13", "25.29", "0.49", "0.0202", 7),
  @(21, "002111", "华宝新起点灵活配置混合", "3.94", "26.43", "0.51", "0.0201", 6),
  @(22, "011224", "九泰盈泰量化股票A", "0.39", "92.77", "4.49", "0.0175", 1),
  @(23, "860029", "光大阳光对冲策略6个月持有期灵活配置混合C", "2.95", "63.41", "0.51", "0.0150", 5),
  @(24, "011225", "九泰盈泰量化股票C", "0.30", "92.77", "4.49", "0.0135", 1),
  @(25, "512530", "建信沪深300红利ETF", "0.48", "98.40", "2.50", "0.0120", 10),
  @(26, "012709", "东方红中证东方红红利低波动指数C", "0.67", "93.80", "1.62", "0.0109", 5),
  @(27, "007939", "华夏网购精选灵活配置混合C", "0.41", "90.71", "2.03", "0.0083", 6),
  @(28, "166402", "浦银安盛沪港深基本面（LOF）", "0.17", "91.35", "3.72", "0.0063", 5),
  @(29, "008093", "同泰慧选混合A", "0.21", "63.44", "2.96", "0.0062", 3),
  @(30, "009208", "建信沪深300指数增强（LOF）C", "0.15", "93.07", "2.85", "0.0043", 4),
  @(31, "008094", "同泰慧选混合C", "0.13", "63.44", "2.96", "0.0038", 3),
  @(32, "002837", "华夏网购精选灵活配置混合A", "0.18", "90.71", "2.03", "0.0037", 6),
  @(33, "162414", "华宝新机遇灵活配置混合（LOF）A", "0.71", "25.29", "0.49", "0.0035", 7),
  @(34, "002334", "汇丰晋信大盘波动精选股票A", "0.16", "85.75", "2.01", "0.0032", 4),
  @(35, "004988", "人保双利优选混合A", "0.56", "25.48", "0.56", "0.0031", 4),
  @(36, "860028", "光大阳光对冲策略6个月持有期灵活配置混合B", "0.39", "63.41", "0.51", "0.0020", 5),
  @(37, "013612", "工银民瑞一年持有混合C", "0.10", "21.97", "1.67", "0.0017", 1),
  @(38, "860010", "光大阳光对冲策略6个月持有期灵活配置混合A", "0.03", "63.41", "0.51", "0.0002", 5),
  @(39, "004989", "人保双利优选混合C", "0.04", "25.48", "0.56", "0.0002", 4),
  @(40, "002335", "汇丰晋信大盘波动精选股票C", "0.01", "85.75", "2.01", "0.0002", 4),
  @(41, "012876", "富荣福耀混合A", "0.00", "61.19", "3.49", "__NUM0__", 8),
)

# ---- 3. First make sure every data row (2..43) in column A carries the
#         same bold/bordered style as the template's row 2, then fill in
#         the 42 fund rows. Columns B..G are text in this workbook
#         (numbers stored with trailing zeros etc.), so we force "@"
#         (Text) format before assigning the string, then clear the
#         left-over number-format style so the cell ends up on the
#         default style (exactly like the source data).
$newSheet.Range("A2").Copy()
$newSheet.Range("A2:A43").PasteSpecial(-4122)

$r = 2
foreach ($row in $rows2022Q3) {
    $idx  = $row[0]
    $code = $row[1]
    $name = $row[2]
    $size = $row[3]
    $pos  = $row[4]
    $pct  = $row[5]
    $mv   = $row[6]
    $rank = $row[7]

    $newSheet.Range("B$r`:F$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $code
    $newSheet.Range("C$r").Value = $name
    $newSheet.Range("D$r").Value = $size
    $newSheet.Range("E$r").Value = $pos
    $newSheet.Range("F$r").Value = $pct
    $newSheet.Range("B$r`:F$r").ClearFormats()

    if ($mv -eq "__NUM0__") {
        $newSheet.Range("G$r").NumberFormat = "General"
        $newSheet.Range("G$r").Value = 0
    } else {
        $newSheet.Range("G$r").NumberFormat = "@"
        $newSheet.Range("G$r").Value = $mv
        $newSheet.Range("G$r").ClearFormats()
    }

    $newSheet.Range("A$r").Value = $idx
    $newSheet.Range("H$r").Value = $rank

    $r++
}


# ---- 4. "总计" (summary) sheet: insert a new row 2 for "2022-Q3" and
#         push the existing quarters down by one. Rows.Insert() leaves
#         B2:D2 carrying a left-over style picked up from the row that
#         used to be there (and no style at all on A2), so fix both up
#         to mirror the sheet's existing convention: column A bold/
#         bordered (style of A3), B..D plain/default.
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 42
$totalSheet.Range("D2").Value = 2.58

Write-Host "2022-Q3 sheet + summary row added."
